# Generate Report for Handback
#
# This script updates the localization-status workbook to reflect a
# "handback" event: the Overview status moves from "Ready for handoff" to
# "Handed back: in sync with en-US", and each per-locale sheet (zh-cn,
# de-de) gets its "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" columns (I, J, K) populated for both data
# rows, including a new hyperlink on column I pointing at the same source
# markdown file as column A. A few columns are also widened to
# accommodate the new, longer datetime-style status text.

$wb = $excel.ActiveWorkbook

$ghBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e759b8715a2075c5232f8cb942e2f11545077889/e2e/"

# ---------------------------------------------------------------------
# Overview sheet: status text + column widths
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

# Columns E (zh-cn) and F (de-de) grow to fit the longer status text.
$wsOverview.Columns.Item(5).ColumnWidth = 29.15
$wsOverview.Columns.Item(6).ColumnWidth = 29.15

# ---------------------------------------------------------------------
# Per-locale sheets: zh-cn (sheet index 2) and de-de (sheet index 3)
# ---------------------------------------------------------------------
$locales = @(
    @{ Sheet = "zh-cn"; Suffix = "zh-cn"; HandbackTime = "2016-09-07 03:00:56" },
    @{ Sheet = "de-de"; Suffix = "de-de"; HandbackTime = "2016-09-07 03:01:12" }
)

$rows = @(
    @{ Row = 2; Uuid = "1dfcb511-9799-41b3-939b-1d76a3bc0009"; Hash = "fa7d65d9b8e3d14f66909ddb11e3da987c52b42e" },
    @{ Row = 3; Uuid = "fc50801d-dc7d-4ccf-bcf8-b1f9258d3670"; Hash = "48d9e223daa15a61650d5c6a0cbc86e876437af5" }
)

foreach ($locale in $locales) {
    $ws = $wb.Worksheets.Item($locale.Sheet)

    # Column C (Status) widens; columns I (Latest Target File) and J
    # (Latest Handback File) widen a lot to fit full file names.
    $ws.Columns.Item(3).ColumnWidth = 29.15
    $ws.Columns.Item(9).ColumnWidth = 39.17
    $ws.Columns.Item(10).ColumnWidth = 39.17

    foreach ($r in $rows) {
        $row = $r.Row
        $mdName = "$($r.Uuid).md"
        $xlfName = "$($r.Uuid).$($r.Hash).$($locale.Suffix).xlf"

        # I<row>: Latest Target File -- same display name/hyperlink as column A.
        $ws.Hyperlinks.Add($ws.Range("I$row"), "$ghBase$mdName", $null, $null, $mdName)

        # J<row>: Latest Handback File -- the generated xliff for this locale.
        $ws.Range("J$row").Value = $xlfName

        # K<row>: Latest Handback DateTime -- replaces the "0001-01-01" placeholder.
        $ws.Range("K$row").Value = $locale.HandbackTime
    }
}
